$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("E1").Value = 'gemini_description'
$ws.Range("F1").Value = 'facial_emotion'
$ws.Range("G1").Value = 'text_similarity_semantic'
$ws.Range("H1").Value = 'text_similarity_semantic_expression'
$ws.Range("I1").Value = 'labels'
$ws.Range("J1").Value = 'bounding_boxes'
$ws.Range("K1").Value = 'bounding_box_confidence'

# Row 2
$ws.Range("B2").Value = 'poster, '
$ws.Range("F2").Value = '---'
$ws.Range("G2").Value = 0.8735178311665853
$ws.Range("H2").Value = 'low'
$ws.Range("J2").Value = '[    0.50049     0.50024     0.98844     0.99073]'
$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = '0.4777811'

# Row 3
$ws.Range("F3").Value = '---'
$ws.Range("G3").Value = 0.8720777829488119
$ws.Range("H3").Value = 'low'

# Row 4
$ws.Range("B4").Value = ""
$ws.Range("F4").Value = '---'
$ws.Range("G4").Value = 0.9031285444895426
$ws.Range("H4").Value = 'medium'

# Row 5
$ws.Range("F5").Value = 'disgust'
$ws.Range("G5").Value = 0.9182257652282715
$ws.Range("H5").Value = 'medium'
$ws.Range("J5").Value = '[    0.55928     0.55452     0.47864     0.57719]'
$ws.Range("K5").NumberFormat = "@"
$ws.Range("K5").Value = '0.5110339'

# Row 6
$ws.Range("B6").Value = 'tree, plant, '
$ws.Range("F6").Value = 'surprise'
$ws.Range("G6").Value = 0.9083539644877116
$ws.Range("H6").Value = 'medium'
$ws.Range("J6").Value = '[     0.4836      0.4533     0.31186     0.37395],[    0.48356     0.45354     0.30985      0.3733]'
$ws.Range("K6").Value = '0.6861728, 0.43398994'

# Row 7
$ws.Range("F7").Value = 'disgust'
$ws.Range("G7").Value = 0.8902324835459391
$ws.Range("H7").Value = 'medium'
$ws.Range("J7").Value = '[    0.64667     0.47426     0.59342     0.71294]'
$ws.Range("K7").NumberFormat = "@"
$ws.Range("K7").Value = '0.51144177'

# Row 8
$ws.Range("B8").Value = 'poster, '
$ws.Range("F8").Value = 'happy'
$ws.Range("G8").Value = 0.9210977554321289
$ws.Range("H8").Value = 'high'
$ws.Range("I8").Value = 'crowded'
$ws.Range("J8").Value = '[    0.25319     0.52039     0.50318     0.94886]'
$ws.Range("K8").NumberFormat = "@"
$ws.Range("K8").Value = '0.57296276'

# Row 9
$ws.Range("B9").Value = 'person, poster, '
$ws.Range("F9").Value = 'surprise'
$ws.Range("G9").Value = 0.8911571502685547
$ws.Range("H9").Value = 'medium'
$ws.Range("J9").Value = '[    0.48131     0.56876     0.58897     0.85699],[    0.42381     0.18348     0.55205     0.36521]'
$ws.Range("K9").Value = '0.49001718, 0.38414583'

# Row 10
$ws.Range("B10").Value = 'tree, '
$ws.Range("F10").Value = '---'
$ws.Range("G10").Value = 0.9032729466756185
$ws.Range("H10").Value = 'medium'
$ws.Range("J10").Value = '[     0.4997     0.75507     0.99399     0.48348]'
$ws.Range("K10").NumberFormat = "@"
$ws.Range("K10").Value = '0.41401386'

# Row 11
$ws.Range("A11").Value = 'id_1104737401602228225_2019-03-10.jpg'
$ws.Range("B11").Value = 'person, '
$ws.Range("C11").Value = 'there is a man standing on a rock in the woods'
$ws.Range("D11").Value = 'A man stands confidently on a rocky cliff, his gaze directed towards the camera. He is dressed in a blue shirt and sunglasses, and the cliff is surrounded by a lush forest of tall trees with green leaves. The man''s position on the cliff and the verdant forest create a sense of depth and perspective in the image.'
$ws.Range("E11").Value = ' This is a photo of a man standing in a forest. The man is in his 40s, with dark hair and brown eyes. He is wearing a blue shirt and sunglasses. He has a friendly expression on his face. The forest is dense, with tall trees and a thick understory. The trees are mostly green, but there are a few yellow and orange leaves. The ground is covered in leaves and moss. There is a small stream in the background. The photo was taken on a sunny day. The background of the photo is a mountain range. The mountains are covered in snow. The photo is taken from a slightly elevated perspective.'
$ws.Range("F11").Value = 'disgust'
$ws.Range("G11").Value = 0.9049077033996582
$ws.Range("H11").Value = 'medium'
$ws.Range("J11").Value = '[    0.21827     0.69374     0.43625     0.61303]'
$ws.Range("K11").NumberFormat = "@"
$ws.Range("K11").Value = '0.6026033'

# Row 12
$ws.Range("A12").Value = 'id_1092817526399078400_2019-02-05.jpg'
$ws.Range("B12").Value = 'poster, person, '
$ws.Range("C12").Value = 'arafed audience of students in a gymnasium watching a man on a stage'
$ws.Range("D12").Value = 'The image depicts a lively scene in a gymnasium, with a group of people gathered in front of a stage. The stage is adorned with a large banner displaying the words "Founded in 1892" and "New York City", suggesting it is a significant event. The gymnasium is filled with people, some sitting on the floor and others standing, all engaged in the event. The perspective of the image is from the back of the gymnasium, providing a comprehensive view of the scene.'
$ws.Range("E12").Value = ' This is a photo of a group of girls in a gym. There are about 100 girls sitting on the bleachers. They are all wearing school uniforms. There is one girl standing in front of the bleachers. She is wearing a black suit. She is holding hands with another girl who is sitting on the bleachers. The girl in the black suit is smiling. The girl sitting on the bleachers is also smiling. The background of the photo is a stage. There is a red curtain behind the stage. There are some chairs on the stage. There is a podium on the stage. There is a banner hanging from the ceiling. The banner has the words "The Young Women''s Leadership School" on it.'
$ws.Range("F12").Value = 'sadness'
$ws.Range("G12").Value = 0.8510868549346924
$ws.Range("H12").Value = 'low'
$ws.Range("I12").Value = 'crowded'
$ws.Range("J12").Value = '[    0.81833      0.3493     0.23254     0.59861],[     0.5434     0.43226     0.17948     0.30161]'
$ws.Range("K12").Value = '0.4040551, 0.4006823'

# Row 13
$ws.Range("A13").Value = 'id_1161027544096923656_2019-08-12.jpg'
$ws.Range("B13").Value = 'person, tree, '
$ws.Range("C13").Value = 'arafed man taking a picture of himself in a garden'
$ws.Range("D13").Value = 'A young man stands in a lush garden, his hands raised to his face as if in deep thought. He is dressed casually in a gray t-shirt and blue shorts. The garden is a vibrant display of nature, with a variety of flowers and plants in shades of green, purple, and white. The man is positioned in the center of the image, surrounded by the verdant foliage.'
$ws.Range("E13").Value = ' This is a photo of a person in a garden. The person is wearing a gray t-shirt, blue and white swim trunks, and glasses. The person has their hands to their ears. There are many plants and flowers in the garden, and a large tree in the background. The background of the image is blurred. The image is taken from a slightly elevated perspective.'
$ws.Range("F13").Value = 'surprise'
$ws.Range("G13").Value = 0.9078138669331869
$ws.Range("H13").Value = 'medium'
$ws.Range("J13").Value = '[    0.21251     0.66714     0.29577     0.31458],[    0.50201     0.41018      0.9868     0.81768]'
$ws.Range("K13").Value = '0.5303586, 0.44902265'

# Row 14
$ws.Range("A14").Value = 'id_1190515268356755461_2019-11-02.jpg'
$ws.Range("B14").Value = 'person, flower, '
$ws.Range("C14").Value = 'there is a man and a woman sitting on a bench'
$ws.Range("D14").Value = 'In the image, there are two individuals seated on a bench in a garden setting. The person on the left is wearing a black baseball cap and a black t-shirt, while the person on the right is wearing a blue jacket and a gray t-shirt. Both individuals are smiling and looking directly at the camera. The garden around them is lush with greenery, including bushes and flowers, and the sky is visible in the background.'
$ws.Range("E14").Value = ' This is a photo of Leonardo DiCaprio and Greta Thunberg. They are sitting outdoors. There are plants in the background. They are both smiling. Leonardo DiCaprio is wearing a black cap and a black t-shirt. Greta Thunberg is wearing a blue sweatshirt and a grey t-shirt.'
$ws.Range("F14").Value = 'disgust'
$ws.Range("G14").Value = 0.9079781373341879
$ws.Range("H14").Value = 'medium'
$ws.Range("J14").Value = '[    0.29657     0.55917      0.5897      0.8774],[    0.95044     0.49172    0.092934    0.093023]'
$ws.Range("K14").Value = '0.4529528, 0.3819592'

# Row 15
$ws.Range("A15").Value = 'id_1113624374702166017_2019-04-04.jpg'
$ws.Range("B15").Value = 'person, '
$ws.Range("C15").Value = 'someone is making a mess in their house with a cardboard box'
$ws.Range("D15").Value = 'In the image, a person is standing in a room, holding a piece of wood and preparing to cut it. The person is wearing a black shirt and is positioned in front of a wooden door. The room has a white wall and a wooden floor, creating a contrast between the person''s dark clothing and the lighter elements of the room. The person is also holding a piece of paper, possibly a message or instructions.'
$ws.Range("E15").Value = ' This is a live video of Alexandria Ocasio-Cortez putting together furniture. She is wearing a black long-sleeved shirt and black pants. She is standing in a room that is mostly empty, except for a few pieces of furniture and some boxes. There is a window in the background. The background of the image is a white wall. The setting of the background is indoor. The type of the image is a video.'
$ws.Range("F15").Value = 'sadness'
$ws.Range("G15").Value = 0.8865760962168375
$ws.Range("H15").Value = 'medium'
$ws.Range("J15").Value = '[    0.17294     0.49255     0.32297     0.43309]'
$ws.Range("K15").NumberFormat = "@"
$ws.Range("K15").Value = '0.4863836'

# Row 16
$ws.Range("A16").Value = 'id_1207450022591238144_2019-12-18.jpg'
$ws.Range("B16").Value = 'car, '
$ws.Range("C16").Value = 'arafed image of a group of men in suits standing next to a race car'
$ws.Range("D16").Value = 'In the image, a group of six individuals dressed in formal attire, including black suits and white shirts, are standing in front of a blue and gray race car. The car is parked in a large warehouse-like space, with a high ceiling and large windows that allow natural light to filter in. The individuals are arranged in a line, with the person in the center standing slightly ahead of the others. The background features a large window that offers a view of the sky, adding depth to the scene.'
$ws.Range("E16").Value = ' There are seven men standing in front of a Formula E race car. The men are all wearing black suits and ties. The car is blue and white. There is a large window in the background. There are also some industrial-looking machines and equipment in the background. The background is dark and shadowy. The image is a photo.'
$ws.Range("F16").Value = 'neutral'
$ws.Range("G16").Value = 0.8991029262542725
$ws.Range("H16").Value = 'medium'
$ws.Range("I16").Value = 'crowded'
$ws.Range("J16").Value = '[    0.32845     0.73155     0.63035     0.25044]'
$ws.Range("K16").NumberFormat = "@"
$ws.Range("K16").Value = '0.49228573'
